$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 172, shifting existing rows 172:182 down to 173:183
$ws.Rows.Item(172).Insert()

# Populate the newly inserted row 172 with the new weekly record
$ws.Range("A172").Value = 1
$ws.Range("B172").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C172").Value = "Arica y Parinacota"
$ws.Range("D172").Value = 44578
$ws.Range("E172").Value = 15
$ws.Range("F172").Value = "Fruta"
$ws.Range("G172").Value = 100108
$ws.Range("H172").Value = "Tropicales y subtropicales"
$ws.Range("I172").Value = 100108006
$ws.Range("J172").Value = "Plátano"
$ws.Range("K172").Value = "Sin especificar"
$ws.Range("L172").Value = "Pintón"
$ws.Range("M172").Value = 120
$ws.Range("N172").Value = 20000
$ws.Range("O172").Value = 22000
$ws.Range("P172").Value = 21000
$ws.Range("Q172").Value = "$/caja 20 kilos"
$ws.Range("R172").Value = "Ecuador"
$ws.Range("S172").Value = 1050
$ws.Range("T172").Value = 20
